# Ran code for averaged intensities on spiral schemes.
#
# This adds three new averaging schemes ("Spiral-90deg-10rot-5space",
# "Spiral-90deg-15rot-5space", "Spiral-90deg-10rot-3space") to the
# alpha4F-HW45 results table, moves "Gaussian-Quadrature" up next to the
# other summary rows, and refreshes the recomputed intensity ratios for
# every averaging-scheme row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the 3 new scheme rows by inserting them right after
#    the existing table (rows 17:19), then extend the worksheet's used
#    range/dimension to A1:M19.
# ---------------------------------------------------------------------
$ws.Rows("17:19").Insert()

# Carry over the bordered/bold "index" formatting used by column A down
# through the newly-inserted rows (matches the style already used by
# A10:A16).
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) New rows 17-19 are the schemes that used to be HexGrid-90degTilt5degRes
#    / HexGrid-90degTilt22p5degRes / HexGrid-60degTilt5degRes (index 15,
#    16, 17) -- they keep their original index + label, just shifted to
#    the bottom of the table.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# ---------------------------------------------------------------------
# 3) Refresh the recomputed averaged-intensity ratios (columns C:M) for
#    every scheme row. Row 10 (Gaussian-Quadrature) now holds the values
#    that used to belong to the old last row; rows 11-13 are the brand
#    new Spiral schemes; rows 14-19 hold the recomputed values for the
#    schemes that used to occupy rows 10-15.
# ---------------------------------------------------------------------

# Row 10 -- Gaussian-Quadrature
$ws.Range("C10").Value = 0.9975163270514803
$ws.Range("D10").Value = 0.9929265524063458
$ws.Range("E10").Value = 0.9988235294117647
$ws.Range("F10").Value = 0.9975163270514803
$ws.Range("G10").Value = 0.9983006600036716
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0.9976470588235294
$ws.Range("J10").Value = 0.9929265524063458
$ws.Range("K10").Value = 0.9958750409090552
$ws.Range("L10").Value = 0.9966956839802676
$ws.Range("M10").Value = 0.9975356879494653

# Row 11 -- Spiral-90deg-10rot-5space
$ws.Range("C11").Value = 0.9965128889843057
$ws.Range("D11").Value = 0.9964745346458684
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.9965128889843057
$ws.Range("G11").Value = 0.9953745366004162
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0.9964745346458684
$ws.Range("K11").Value = 0.9982372673229343
$ws.Range("L11").Value = 0.9973750781536199
$ws.Range("M11").Value = 0.9980603267050983

# Row 12 -- Spiral-90deg-15rot-5space
$ws.Range("C12").Value = 0.9965062734929344
$ws.Range("D12").Value = 0.9965447663052212
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.9965062734929344
$ws.Range("G12").Value = 0.995356727983477
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0.9965447663052212
$ws.Range("K12").Value = 0.9982723831526106
$ws.Range("L12").Value = 0.9973893283227724
$ws.Range("M12").Value = 0.9980679612969388

# Row 13 -- Spiral-90deg-10rot-3space
$ws.Range("C13").Value = 0.9965140721589814
$ws.Range("D13").Value = 0.9964890192639335
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.9965140721589814
$ws.Range("G13").Value = 0.9953656277827446
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 0.9964890192639335
$ws.Range("K13").Value = 0.9982445096319668
$ws.Range("L13").Value = 0.9973792908954741
$ws.Range("M13").Value = 0.9980614532009433

# Row 14 -- NoRotation-tilt60deg
$ws.Range("C14").Value = 0.9975039999999989
$ws.Range("D14").Value = 0.9883280000000004
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.9975039999999989
$ws.Range("G14").Value = 0.9983359999999993
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.9883280000000004
$ws.Range("K14").Value = 0.9941640000000003
$ws.Range("L14").Value = 0.9958339999999996
$ws.Range("M14").Value = 0.9973613333333331

# Row 15 -- Rotation-NoTilt
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.98
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0.98
$ws.Range("K15").Value = 0.99
$ws.Range("L15").Value = 0.995
$ws.Range("M15").Value = 0.9966666666666667

# Row 16 -- Rotation-60detTilt
$ws.Range("C16").Value = 0.9992324802560002
$ws.Range("D16").Value = 0.9882299580416009
$ws.Range("E16").Value = 0.9986692104191994
$ws.Range("F16").Value = 0.9992324802560002
$ws.Range("G16").Value = 0.9981859385344
$ws.Range("H16").Value = 0.9993626025984018
$ws.Range("I16").Value = 0.9975945312255989
$ws.Range("J16").Value = 0.9882299580416009
$ws.Range("K16").Value = 0.9934495842304001
$ws.Range("L16").Value = 0.9963410322432003
$ws.Range("M16").Value = 0.9968791201792002

# Row 17 -- HexGrid-90degTilt5degRes
$ws.Range("C17").Value = 0.9981472497059053
$ws.Range("D17").Value = 0.9989293515472001
$ws.Range("E17").Value = 0.9962099393971589
$ws.Range("F17").Value = 0.9981472497059053
$ws.Range("G17").Value = 0.9967250324794205
$ws.Range("H17").Value = 0.9985870440767033
$ws.Range("I17").Value = 0.9937878340704579
$ws.Range("J17").Value = 0.9989293515472001
$ws.Range("K17").Value = 0.9975696454721794
$ws.Range("L17").Value = 0.9978584475890424
$ws.Range("M17").Value = 0.9970644085461412

# Row 18 -- HexGrid-90degTilt22p5degRes
$ws.Range("C18").Value = 0.9979231908498345
$ws.Range("D18").Value = 0.9986046511627907
$ws.Range("E18").Value = 0.9947059819900204
$ws.Range("F18").Value = 0.9979231908498345
$ws.Range("G18").Value = 0.9964139249229853
$ws.Range("H18").Value = 0.9988381213999441
$ws.Range("I18").Value = 0.9935760165628401
$ws.Range("J18").Value = 0.9986046511627907
$ws.Range("K18").Value = 0.9966553165764056
$ws.Range("L18").Value = 0.9972892537131202
$ws.Range("M18").Value = 0.9966769811480692

# Row 19 -- HexGrid-60degTilt5degRes
$ws.Range("C19").Value = 0.9976509728865289
$ws.Range("D19").Value = 0.9992225252278324
$ws.Range("E19").Value = 0.9965128213365479
$ws.Range("F19").Value = 0.9976509728865289
$ws.Range("G19").Value = 0.9977770030111814
$ws.Range("H19").Value = 0.9978926920620124
$ws.Range("I19").Value = 0.9936302288237914
$ws.Range("J19").Value = 0.9992225252278324
$ws.Range("K19").Value = 0.9978676732821902
$ws.Range("L19").Value = 0.9977593230843596
$ws.Range("M19").Value = 0.9971143738913159
